$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the element description text (shared string referenced by A2):
# "Ladrillo rojo de  3 x3 " -> "Tornillos AA"
$ws.Range("A2").Value = "Tornillos AA"

# Update numeric values in row 2
$ws.Range("C2").Value = 111
$ws.Range("D2").Value = 122

# Remove the value previously stored in S2
[void]$ws.Range("S2").ClearContents()

# Update the current selection to C2 (also drops the stale topLeftCell scroll anchor)
[void]$ws.Range("C2").Select()

Write-Output "done"
